$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 55, shifting existing rows 55:82 down to 56:83.
$ws.Rows.Item(55).Insert()

# The inserted row keeps the data that was in the old row 55 (Excel's
# Insert copies formatting/values are NOT duplicated by default - cells
# start blank), so re-populate it explicitly matching the original row's
# unchanged fields and the new record's changed fields.
$ws.Cells.Item(55, 1).Value = 5
$ws.Cells.Item(55, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(55, 3).Value = "Maule"
$ws.Cells.Item(55, 4).Value = 44875
$ws.Cells.Item(55, 5).Value = 7
$ws.Cells.Item(55, 6).Value = 300000000
$ws.Cells.Item(55, 7).Value = "Espárragos"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 3000
$ws.Cells.Item(55, 11).Value = 1000
$ws.Cells.Item(55, 12).Value = 1000
$ws.Cells.Item(55, 13).Value = 1000
$ws.Cells.Item(55, 14).Value = "$/kilo"
$ws.Cells.Item(55, 15).Value = "Provincia de Linares"
$ws.Cells.Item(55, 16).Value = 1000
$ws.Cells.Item(55, 17).Value = 1
$ws.Cells.Item(55, 18).Value = "Hortaliza"
